$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.7265917602996255
$ws1.Range("C2").Value = 0.9583333333333334
$ws1.Range("D2").Value = 0.4737827715355805
$ws1.Range("E2").Value = 0.6340852130325815
$ws1.Range("F2").Value = 0.5270833333333333
$ws1.Range("G2").Value = 0.4831790803584545
$ws1.Range("H2").Value = 0.7265917602996254
$ws1.Range("I2").Value = 253
$ws1.Range("J2").Value = 11
$ws1.Range("K2").Value = 523
$ws1.Range("L2").Value = 281

# --- Sheet: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")
$ws2.Range("B2").Value = 0.650497512437811
$ws2.Range("C2").Value = 0.9794007490636704
$ws2.Range("D2").Value = 0.7817638266068759

$ws2.Range("B3").Value = 0.9583333333333334
$ws2.Range("C3").Value = 0.4737827715355805
$ws2.Range("D3").Value = 0.6340852130325815

$ws2.Range("B4").Value = 0.7265917602996255
$ws2.Range("C4").Value = 0.7265917602996255
$ws2.Range("D4").Value = 0.7265917602996255
$ws2.Range("E4").Value = 0.7265917602996255

$ws2.Range("B5").Value = 0.8044154228855722
$ws2.Range("C5").Value = 0.7265917602996255
$ws2.Range("D5").Value = 0.7079245198197287

$ws2.Range("B6").Value = 0.8044154228855722
$ws2.Range("C6").Value = 0.7265917602996255
$ws2.Range("D6").Value = 0.7079245198197286

# --- Sheet: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 523
$ws3.Range("C2").Value = 11
$ws3.Range("B3").Value = 281
$ws3.Range("C3").Value = 253
